$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix "Trusted by Leading Industry Professionals" -> "Acknowledged by the industry’s leading professionals" (row 14 / thirdSectionHeader)
$ws.Range("B14").Value = "Acknowledged by the industry’s leading professionals"

# Fix "Lifetime" -> "Yearly" button text (row 4 / firstSectionbutton2)
$ws.Range("B4").Value = "Yearly"

# Update the active selection to B15
$ws.Range("B15").Select() | Out-Null
